$wb = $excel.ActiveWorkbook

# --- Worksheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 602.7143
$ws.Range("I33").Value = 639.9231
$ws.Range("J33").Value = 119
$ws.Range("K33").Value = 639.9231
$ws.Range("L33").Value = 119
$ws.Range("M33").Value = -410.9231
$ws.Range("N33").Value = -577
$ws.Range("H98").Value = 1293.0476
$ws.Range("I98").Value = 1362.5883
$ws.Range("J98").Value = 997.5
$ws.Range("K98").Value = 1362.5883
$ws.Range("L98").Value = 997.5
$ws.Range("M98").Value = 135.4117000000001
$ws.Range("N98").Value = -3993.5
$ws.Range("H122").Value = 1293.0476
$ws.Range("I122").Value = 1362.5883
$ws.Range("J122").Value = 997.5
$ws.Range("K122").Value = 4087.7649
$ws.Range("L122").Value = 2992.5
$ws.Range("M122").Value = -1637.7649
$ws.Range("N122").Value = -7892.5
$ws.Range("H137").Value = 1315.1025
$ws.Range("I137").Value = 1158.3549
$ws.Range("K137").Value = 3475.0647
$ws.Range("M137").Value = -925.0646999999999
$ws.Range("H138").Value = 1656.262
$ws.Range("I138").Value = 1090.1578
$ws.Range("J138").Value = 2123.913
$ws.Range("K138").Value = 3270.4734
$ws.Range("L138").Value = 6371.739
$ws.Range("M138").Value = 1869.5266
$ws.Range("N138").Value = -16651.739
$ws.Range("H141").Value = 4917.1924
$ws.Range("I141").Value = 1542.6923
$ws.Range("J141").Value = 8291.691999999999
$ws.Range("K141").Value = 4628.0769
$ws.Range("L141").Value = 24875.076
$ws.Range("M141").Value = 551.9231
$ws.Range("N141").Value = -35235.076

# --- Worksheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1178.2593
$ws.Range("I2").Value = 918.75
$ws.Range("J2").Value = 1555.7273
$ws.Range("K2").Value = 918.75
$ws.Range("L2").Value = 1555.7273
$ws.Range("M2").Value = -805.75
$ws.Range("N2").Value = -1781.7273
$ws.Range("H74").Value = 2078.318
$ws.Range("I74").Value = 1766.0588
$ws.Range("K74").Value = 1766.0588
$ws.Range("M74").Value = -892.0588
$ws.Range("H77").Value = 2078.318
$ws.Range("I77").Value = 1766.0588
$ws.Range("K77").Value = 8830.294
$ws.Range("M77").Value = -4462.294
$ws.Range("H116").Value = 1178.2593
$ws.Range("I116").Value = 918.75
$ws.Range("J116").Value = 1555.7273
$ws.Range("K116").Value = 918.75
$ws.Range("L116").Value = 1555.7273
$ws.Range("M116").Value = 1375.25
$ws.Range("N116").Value = -6143.7273

# --- Worksheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1178.2593
$ws.Range("I3").Value = 918.75
$ws.Range("J3").Value = 1555.7273
$ws.Range("K3").Value = 918.75
$ws.Range("L3").Value = 1555.7273
$ws.Range("M3").Value = -804.75
$ws.Range("N3").Value = -1783.7273
$ws.Range("H110").Value = 32857.145
$ws.Range("J110").Value = 32857.145
$ws.Range("L110").Value = 32857.145
$ws.Range("N110").Value = -41037.145
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -33134

# --- Worksheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3848166.5
$ws.Range("I31").Value = 2003.4546
$ws.Range("J31").Value = 25002062
$ws.Range("K31").Value = 2003.4546
$ws.Range("L31").Value = 25002062
$ws.Range("M31").Value = -1708.4546
$ws.Range("N31").Value = -25002652
$ws.Range("H34").Value = 3848166.5
$ws.Range("I34").Value = 2003.4546
$ws.Range("J34").Value = 25002062
$ws.Range("K34").Value = 2003.4546
$ws.Range("L34").Value = 25002062
$ws.Range("M34").Value = -1801.4546
$ws.Range("N34").Value = -25002466
$ws.Range("H58").Value = 1209.2174
$ws.Range("I58").Value = 1167.6666
$ws.Range("J58").Value = 1254.5454
$ws.Range("K58").Value = 1167.6666
$ws.Range("L58").Value = 1254.5454
$ws.Range("M58").Value = -964.6666
$ws.Range("N58").Value = -1660.5454
$ws.Range("H86").Value = 2235.375
$ws.Range("I86").Value = 2101.75
$ws.Range("K86").Value = 2101.75
$ws.Range("M86").Value = -978.75
$ws.Range("H89").Value = 2235.375
$ws.Range("I89").Value = 2101.75
$ws.Range("K89").Value = 10508.75
$ws.Range("M89").Value = -4892.75
$ws.Range("H134").Value = 915.35
$ws.Range("I134").Value = 891.17645
$ws.Range("J134").Value = 1052.3334
$ws.Range("K134").Value = 2673.52935
$ws.Range("L134").Value = 3157.0002
$ws.Range("M134").Value = -138.5293500000002
$ws.Range("N134").Value = -8227.0002
$ws.Range("H136").Value = 1209.2174
$ws.Range("I136").Value = 1167.6666
$ws.Range("J136").Value = 1254.5454
$ws.Range("K136").Value = 3502.9998
$ws.Range("L136").Value = 3763.6362
$ws.Range("M136").Value = -952.9998000000001
$ws.Range("N136").Value = -8863.636200000001

# --- Worksheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1012246.25
$ws.Range("J122").Value = 1390325.8
$ws.Range("L122").Value = 12512932.2
$ws.Range("N122").Value = -12517832.2
$ws.Range("H131").Value = 3178119
$ws.Range("I131").Value = 6682.3125
$ws.Range("J131").Value = 5848802.5
$ws.Range("K131").Value = 20046.9375
$ws.Range("L131").Value = 17546407.5
$ws.Range("M131").Value = -15006.9375
$ws.Range("N131").Value = -17556487.5

# --- Worksheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3768750
$ws.Range("I7").Value = 3768750
$ws.Range("K7").Value = 3768750
$ws.Range("M7").Value = -3768638
$ws.Range("H8").Value = 3768750
$ws.Range("I8").Value = 3768750
$ws.Range("K8").Value = 3768750
$ws.Range("M8").Value = -3768611
$ws.Range("H126").Value = 8334836
$ws.Range("I126").Value = 1556
$ws.Range("J126").Value = 16668116
$ws.Range("K126").Value = 4668
$ws.Range("L126").Value = 50004348
$ws.Range("M126").Value = -2198
$ws.Range("N126").Value = -50009288

# --- Worksheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 9630
$ws.Range("J109").Value = 9630
$ws.Range("L109").Value = 9630
$ws.Range("N109").Value = -12404
$ws.Range("H133").Value = 24734
$ws.Range("J133").Value = 24734
$ws.Range("L133").Value = 24734
$ws.Range("N133").Value = -29794
$ws.Range("H136").Value = 5932
$ws.Range("I136").Value = 7281.5
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 21844.5
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = -19294.5
$ws.Range("N136").Value = -12099.9999
